$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Interrupt every" section using a 50Hz RAD timer instead of BPM-derived
# timing. Write the shared-string labels in the same order the original
# author entered them (F4, F1, F2) so the sharedStrings table lines up.
$ws.Range("F4").Value = "Interrupt every"

$ws.Range("F1").Value = "RAD Timer (Hz)"
$ws.Range("G1").Value = 50
$ws.Range("H1").Formula = "=1/G1"

$ws.Range("F2").Value = "CPU Clock (Hz)"
$ws.Range("G2").Value = 14318000

$ws.Range("H4").Formula = "=H1 * G2"

# Best-effort cosmetic tweak matching the author's wider column F.
$ws.Columns.Item(6).ColumnWidth = 13.8

# Match the author's final cursor position/selection.
[void]$ws.Range("I11").Select()
